# unify the conception of DataNode, DataTable, Entity.
# Rename the two worksheets to reflect the new naming scheme and leave the
# second sheet (DataTable) as the active/selected sheet, matching the
# state the workbook was saved in.

$wb = $excel.ActiveWorkbook

$wsNode  = $wb.Worksheets.Item(1)   # was "Property1"
$wsTable = $wb.Worksheets.Item(2)   # was "Record"

$wsNode.Name  = "DataNode"
$wsTable.Name = "DataTable"

# The workbook was re-saved with the DataTable sheet selected/active.
$wsTable.Activate()
